# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Row -> new F-column value (applies to both the 展览 and 全部类型 sheets).
$updates = @{
    2  = 11809
    3  = 11532
    6  = 1040
    8  = 72
    11 = 10848
    13 = 20
    16 = 2473
    17 = 1056
    20 = 134
    22 = 11157
    23 = 10953
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
